$wb = $excel.ActiveWorkbook

# The "Turkey" sheet is the last sheet and has the same layout (columns A, B, D
# only, with bestFit widths, merged header cells, etc.) that the new "Croatia"
# sheet uses. Duplicate it and place the copy right after "Turkey" so it
# becomes the new last (active) tab, exactly like the new Croatia worksheet
# that was added to the workbook.
$turkey = $wb.Worksheets.Item("Turkey")
$turkey.Copy($null, $turkey)

# The copy is inserted immediately after Turkey - grab it and rename it.
$croatia = $wb.Worksheets.Item($turkey.Index + 1)
$croatia.Name = "Croatia"

# Fill in the market-specific values. B4 (the Jira/NGC reference) is set
# before B2 (the market name) so that any newly created shared strings are
# appended to xl/sharedStrings.xml in the same order as the target workbook
# ("NGC-3139/T2478" before "Croatia Market").
$croatia.Range("B4").Value = "NGC-3139/T2478"
$croatia.Range("B2").Value = "Croatia Market"

# Match the active selection on the new sheet.
$croatia.Range("B2").Select()
